$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.255.12'
$ws.Range("E2").Value = '  -0.10%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.613.87'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '594.36'
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.00'
$ws.Range("E6").Value = '  -0.93%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.553'
$ws.Range("E8").Value = '  +1.18%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.608.93'
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -3.10%  '
$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.20'
$ws.Range("E12").Value = '  -1.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.346'
$ws.Range("E13").Value = '  -2.73%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.66'
$ws.Range("E14").Value = '  -0.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.088.16'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000180'
$ws.Range("E16").Value = '  -3.62%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '66.973.36'
$ws.Range("E17").Value = '  -0.88%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.610.37'
$ws.Range("E18").Value = '  -0.08%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '366.88'
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.07'
$ws.Range("E20").Value = '  -2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.39'
$ws.Range("E21").Value = '  -4.21%  '
$ws.Range("E22").Value = '  -0.30%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.07'
$ws.Range("E23").Value = '  +2.68%  '
$ws.Range("E24").Value = '  -0.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '10.14'
$ws.Range("E25").Value = '  +2.39%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '67.71'
$ws.Range("E26").Value = '  -0.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.739.43'
$ws.Range("E27").Value = '  -0.09%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '587.99'
$ws.Range("E28").Value = '  +1.92%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0000101'
$ws.Range("E30").Value = '  -2.90%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.39'
$ws.Range("E31").Value = '  -3.17%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.74'
$ws.Range("E32").Value = '  -2.71%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.81'
$ws.Range("E33").Value = '  -2.40%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.09%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.123'
$ws.Range("E35").Value = '  -6.91%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("E36").Value = '  -1.42%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '4.85'
$ws.Range("E37").Value = '  -1.38%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '155.51'
$ws.Range("E38").Value = '  +0.11%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.94'
$ws.Range("E39").Value = '  -2.46%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.366'
$ws.Range("E40").Value = '  -0.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.25'
$ws.Range("E41").Value = '  -2.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.81'
$ws.Range("E42").Value = '  -2.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.60'
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.70'
$ws.Range("E44").Value = '  +1.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.86'
$ws.Range("E45").Value = '  -1.12%  '
$ws.Range("E46").Value = '  +0.00%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '154.67'
$ws.Range("E47").Value = '  -1.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0294'
$ws.Range("E48").Value = '  +1.32%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.72'
$ws.Range("E49").Value = '  -0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '21.67'
$ws.Range("E50").Value = '  +3.74%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.617'
$ws.Range("E51").Value = '  -1.92%  '
